# Getting table data from API (Rohan changes)
# The Table component on row 12 switches from static columns/rows JSON
# props to an "api" prop pointing at a REST endpoint, keeping the
# stickyHeader prop but moving it left, and dropping the "mode" column J.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: Table component props
$ws.Range("E12").Value = "api"
$ws.Range("F12").Value = "http://localhost:9001/tableData"
$ws.Range("G12").Value = "stickyHeader"
$ws.Range("H12").Value = $true
$ws.Range("I12").ClearContents()
$ws.Range("J12").ClearContents()

# Update the active selection to match the saved workbook state
$ws.Range("G15").Select()
